$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.331717967987061
$ws.Range("B1").Value = 5.509922981262207
$ws.Range("C1").Value = 4.310285091400146
$ws.Range("D1").Value = 1.879854083061218
$ws.Range("E1").Value = 1.298789858818054
